$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: new pairing test case with no related file info
$ws.Range("B13").Value = "test_lab:f5_1, test_lab:alt_f5_1"
$ws.Range("C13").Value = "fastq"
$ws.Range("D13").Value = ""
$ws.Range("E13").Value = ""

# Row 14: fully blank row (styled only)
$ws.Range("B14").Value = ""
$ws.Range("C14").Value = ""
$ws.Range("D14").Value = ""
$ws.Range("E14").Value = ""

# Copy the style of an existing "text label" cell (B12, style s="2") down to B13/B14
$ws.Range("B12").Copy() | Out-Null
$ws.Range("B13:B14").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Copy the style of an existing "value" cell (C12:E12, style s="1") down to C13:E14
$ws.Range("C12:E12").Copy() | Out-Null
$ws.Range("C13:E14").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$excel.CutCopyMode = 0

$ws.Range("C27").Select() | Out-Null
